# edit.ps1
# Applies the changes described by the commit diff:
#   1. Rename worksheet "My Series" -> "Data"
#   2. Change custom number format 165 (0.000) to ###0.000 for the
#      numeric data column (B13:B29) that used it.
#   3. Change cell A11 text "Function Description" -> "Function Information"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet
$ws.Name = "Data"

# 2. Update the number format applied to the data column (previously "0.000")
$ws.Range("B13:B29").NumberFormat = "###0.000"

# 3. Update the label cell text
$ws.Range("A11").Value = "Function Information"
